# Update scripts with new TPM-derived NATMI metrics for Fn1-Cd44 LR pair.
# Columns: G=Ligand avg expr, H=Ligand total expr, I/J=Ligand specificity (avg/total),
#          M=Receptor avg expr, N=Receptor total expr, O/P=Receptor specificity (avg/total),
#          Q=Edge avg weight, R=Edge total weight, S/T=Edge specificity (avg/total)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 7).Value = 6.240107999999999
$ws.Cells.Item(2, 8).Value = 18.720324
$ws.Cells.Item(2, 9).Value = 0.01732230523539376
$ws.Cells.Item(2, 10).Value = 0.01732230523539376
$ws.Cells.Item(2, 13).Value = 8.142376000000001
$ws.Cells.Item(2, 14).Value = 24.427128
$ws.Cells.Item(2, 15).Value = 0.1741313933276368
$ws.Cells.Item(2, 16).Value = 0.1741313933276368
$ws.Cells.Item(2, 17).Value = 50.809305616608
$ws.Cells.Item(2, 18).Value = 457.283750549472
$ws.Cells.Item(2, 19).Value = 0.003016357146285733
$ws.Cells.Item(2, 20).Value = 0.003016357146285733
$ws.Cells.Item(3, 7).Value = 6.240107999999999
$ws.Cells.Item(3, 8).Value = 18.720324
$ws.Cells.Item(3, 9).Value = 0.01732230523539376
$ws.Cells.Item(3, 10).Value = 0.01732230523539376
$ws.Cells.Item(3, 15).Value = 0.5205382400466131
$ws.Cells.Item(3, 16).Value = 0.5205382400466131
$ws.Cells.Item(3, 17).Value = 151.886377397188
$ws.Cells.Item(3, 18).Value = 1366.977396574692
$ws.Cells.Item(3, 19).Value = 0.009016922280782099
$ws.Cells.Item(3, 20).Value = 0.009016922280782101
$ws.Cells.Item(4, 7).Value = 6.240107999999999
$ws.Cells.Item(4, 8).Value = 18.720324
$ws.Cells.Item(4, 9).Value = 0.01732230523539376
$ws.Cells.Item(4, 10).Value = 0.01732230523539376
$ws.Cells.Item(4, 15).Value = 0.3053303666257501
$ws.Cells.Item(4, 16).Value = 0.3053303666257501
$ws.Cells.Item(4, 17).Value = 89.09148210127199
$ws.Cells.Item(4, 18).Value = 801.8233389114479
$ws.Cells.Item(4, 19).Value = 0.005289025808325928
$ws.Cells.Item(4, 20).Value = 0.005289025808325929
$ws.Cells.Item(5, 9).Value = 0.9592798330716089
$ws.Cells.Item(5, 10).Value = 0.9592798330716091
$ws.Cells.Item(5, 13).Value = 8.142376000000001
$ws.Cells.Item(5, 14).Value = 24.427128
$ws.Cells.Item(5, 15).Value = 0.1741313933276368
$ws.Cells.Item(5, 16).Value = 0.1741313933276368
$ws.Cells.Item(5, 17).Value = 2813.733019251704
$ws.Cells.Item(5, 18).Value = 25323.59717326534
$ws.Cells.Item(5, 19).Value = 0.1670407339238621
$ws.Cells.Item(5, 20).Value = 0.1670407339238621
$ws.Cells.Item(6, 9).Value = 0.9592798330716089
$ws.Cells.Item(6, 10).Value = 0.9592798330716091
$ws.Cells.Item(6, 15).Value = 0.5205382400466131
$ws.Cells.Item(6, 16).Value = 0.5205382400466131
$ws.Cells.Item(6, 19).Value = 0.4993418360193041
$ws.Cells.Item(6, 20).Value = 0.4993418360193042
$ws.Cells.Item(7, 9).Value = 0.9592798330716089
$ws.Cells.Item(7, 10).Value = 0.9592798330716091
$ws.Cells.Item(7, 15).Value = 0.3053303666257501
$ws.Cells.Item(7, 16).Value = 0.3053303666257501
$ws.Cells.Item(7, 19).Value = 0.2928972631284427
$ws.Cells.Item(7, 20).Value = 0.2928972631284428
$ws.Cells.Item(8, 7).Value = 8.428738666666666
$ws.Cells.Item(8, 9).Value = 0.02339786169299727
$ws.Cells.Item(8, 10).Value = 0.02339786169299728
$ws.Cells.Item(8, 13).Value = 8.142376000000001
$ws.Cells.Item(8, 14).Value = 24.427128
$ws.Cells.Item(8, 15).Value = 0.1741313933276368
$ws.Cells.Item(8, 16).Value = 0.1741313933276368
$ws.Cells.Item(8, 17).Value = 68.62995942973866
$ws.Cells.Item(8, 18).Value = 617.669634867648
$ws.Cells.Item(8, 19).Value = 0.004074302257488953
$ws.Cells.Item(8, 20).Value = 0.004074302257488954
$ws.Cells.Item(9, 7).Value = 8.428738666666666
$ws.Cells.Item(9, 9).Value = 0.02339786169299727
$ws.Cells.Item(9, 10).Value = 0.02339786169299728
$ws.Cells.Item(9, 15).Value = 0.5205382400466131
$ws.Cells.Item(9, 16).Value = 0.5205382400466131
$ws.Cells.Item(9, 19).Value = 0.01217948174652687
$ws.Cells.Item(9, 20).Value = 0.01217948174652687
$ws.Cells.Item(10, 7).Value = 8.428738666666666
$ws.Cells.Item(10, 9).Value = 0.02339786169299727
$ws.Cells.Item(10, 10).Value = 0.02339786169299728
$ws.Cells.Item(10, 15).Value = 0.3053303666257501
$ws.Cells.Item(10, 16).Value = 0.3053303666257501
$ws.Cells.Item(10, 17).Value = 89.09148210127199
$ws.Cells.Item(10, 18).Value = 801.8233389114479
$ws.Cells.Item(10, 19).Value = 0.007144077688981451
$ws.Cells.Item(10, 20).Value = 0.007144077688981453
